# Applies the "Improved Excel Export, added FileDialogWithLastOpen class" edit
# to the celldatabaseV2 excel template:
#  - adds a "RpDark[ohm]" column header on the result sheet (G1)
#  - adds a "#UIChartLightDark" tag cell (E6)
#  - relocates the "#UIChart"/"#UIChartDark" tag cells from row 10 (D10/I10)
#    down to a new row 33 (D33/H33)
#  - updates the used range / selection bookkeeping to match

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("templateResultSheet")

# G1 was a duplicate of F1 ("Rp[ohm]"); it should read "RpDark[ohm]" to line
# up with the "#rpDark" tag already present in G2.
$ws1.Range("G1").Value2 = "RpDark[ohm]"

# New tag cell next to the dark-current voltage/current block.
$ws1.Range("E6").Value2 = "#UIChartLightDark"

# Move the chart tag cells from row 10 down to row 33.
$ws1.Range("D10").ClearContents()
$ws1.Range("I10").ClearContents()
$ws1.Range("D33").Value2 = "#UIChart"
$ws1.Range("H33").Value2 = "#UIChartDark"

# Match the author's final selection on the result sheet.
$ws1.Range("D33:H33").Select()
